# Réseau V1 (DMVPN, VLANs, RipV2)
#
# - Corrects the "Responsables association" subnet address in E4
#   (was a typo: 172.2.0.0 -> 10.2.0.0)
# - Adds a new "Routeurs VPN" subnet as row 5 (/25 = 255.255.255.128,
#   range 10.2.4.0 - 10.2.4.127)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the E4 value (typo 172.2.0.0 -> 10.2.0.0)
$ws.Range("E4").Value = "10.2.0.0"

# Clone row 4's formatting (borders/fonts/number formats) into the new row 5
$ws.Range("A4:H4").Copy()
$ws.Range("A5:H5").PasteSpecial(-4122)

# Populate the new "Routeurs VPN" row
$ws.Range("D5").Value = "/25 = 255.255.255.128"
$ws.Range("E5").Value = "10.2.4.0"
$ws.Range("F5").Value = "10.2.4.1"
$ws.Range("G5").Value = "10.2.4.126"
$ws.Range("H5").Value = "10.2.4.127"
$ws.Range("B5").Value = 100
$ws.Range("C5").Value = 126
$ws.Range("A5").Value = "Routeurs VPN"

$ws.Rows.Item(5).RowHeight = 15.75

$ws.Range("E5").Select()
